# Update column G ("K" - strikeouts) values on Sheet1 to reflect the
# regenerated save_data (K instead of Strike#).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 0
    6  = 1
    7  = 2
    8  = 0
    9  = 1
    10 = 2
    11 = 2
    12 = 2
    13 = 4
    14 = 5
    15 = 2
    16 = 3
    17 = 4
    18 = 5
    19 = 5
    20 = 1
    21 = 1
    22 = 1
    23 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
